$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 6868.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 6868.5
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 20605.5
$ws.Range("N69").Value = -22353.5
$ws.Range("M69").ClearContents()

$ws.Range("H70").Value = 4203.1875
$ws.Range("I70").Value = 4795.5454
$ws.Range("J70").Value = 2900
$ws.Range("K70").Value = 14386.6362
$ws.Range("L70").Value = 8700
$ws.Range("M70").Value = -14116.6362
$ws.Range("N70").Value = -9240

$ws.Range("H72").Value = 6868.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 6868.5
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 61816.5
$ws.Range("N72").Value = -70552.5
$ws.Range("M72").ClearContents()

$ws.Range("H73").Value = 4203.1875
$ws.Range("I73").Value = 4795.5454
$ws.Range("J73").Value = 2900
$ws.Range("K73").Value = 14386.6362
$ws.Range("L73").Value = 8700
$ws.Range("M73").Value = -13450.6362
$ws.Range("N73").Value = -10572

$ws.Range("H107").Value = 896.45
$ws.Range("I107").Value = 944.7778
$ws.Range("K107").Value = 944.7778
$ws.Range("M107").Value = 975.2222

$ws.Range("H135").Value = 1217.2963
$ws.Range("I135").Value = 1025.7916
$ws.Range("K135").Value = 9232.124400000001
$ws.Range("M135").Value = -6697.124400000001

$ws.Range("H138").Value = 8199796.5
$ws.Range("J138").Value = 10420140
$ws.Range("L138").Value = 31260420
$ws.Range("N138").Value = -31270700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 8539.799999999999
$ws.Range("I10").Value = 350
$ws.Range("J10").Value = 13999.667
$ws.Range("K10").Value = 350
$ws.Range("L10").Value = 13999.667
$ws.Range("M10").Value = -180
$ws.Range("N10").Value = -14339.667

$ws.Range("H32").Value = 6547.7954
$ws.Range("I32").Value = 3271.1282
$ws.Range("K32").Value = 3271.1282
$ws.Range("M32").Value = -2984.1282

$ws.Range("H33").Value = 15500
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H45").Value = 5374.654
$ws.Range("I45").Value = 6449.75
$ws.Range("K45").Value = 6449.75
$ws.Range("M45").Value = -6072.75

$ws.Range("H61").Value = 221493.92
$ws.Range("I61").Value = 3289.6365
$ws.Range("K61").Value = 3289.6365
$ws.Range("M61").Value = -3077.6365

$ws.Range("H102").Value = 2244.2856
$ws.Range("I102").Value = 2276.55
$ws.Range("K102").Value = 2276.55
$ws.Range("M102").Value = -654.5500000000002

$ws.Range("H110").Value = 5091.756
$ws.Range("I110").Value = 5156.4243
$ws.Range("J110").Value = 4825
$ws.Range("K110").Value = 5156.4243
$ws.Range("L110").Value = 4825
$ws.Range("M110").Value = -3111.4243
$ws.Range("N110").Value = -8915

$ws.Range("H136").Value = 221493.92
$ws.Range("I136").Value = 3289.6365
$ws.Range("K136").Value = 9868.9095
$ws.Range("M136").Value = -7318.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 5250
$ws.Range("J23").Value = 5250
$ws.Range("L23").Value = 5250
$ws.Range("N23").Value = -5816

$ws.Range("H58").Value = 27644.75
$ws.Range("J58").Value = 27644.75
$ws.Range("L58").Value = 27644.75
$ws.Range("N58").Value = -28232.75

$ws.Range("H64").Value = 1935.4286
$ws.Range("J64").Value = 1935.4286
$ws.Range("L64").Value = 1935.4286
$ws.Range("N64").Value = -2385.4286

$ws.Range("H67").Value = 1935.4286
$ws.Range("J67").Value = 1935.4286
$ws.Range("L67").Value = 1935.4286
$ws.Range("N67").Value = -3495.4286

$ws.Range("H94").Value = 859.225
$ws.Range("I94").Value = 530.05884
$ws.Range("J94").Value = 2724.5
$ws.Range("K94").Value = 530.05884
$ws.Range("L94").Value = 2724.5
$ws.Range("M94").Value = -79.05884000000003
$ws.Range("N94").Value = -3626.5

$ws.Range("H134").Value = 1715.9143
$ws.Range("I134").Value = 1538.3939
$ws.Range("K134").Value = 4615.1817
$ws.Range("M134").Value = -2080.1817

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1072.7858
$ws.Range("I5").Value = 604.75
$ws.Range("J5").Value = 1260
$ws.Range("K5").Value = 604.75
$ws.Range("L5").Value = 1260
$ws.Range("M5").Value = -492.75
$ws.Range("N5").Value = -1484

$ws.Range("H28").Value = 30487
$ws.Range("J28").Value = 30487
$ws.Range("L28").Value = 30487
$ws.Range("N28").Value = -30977

$ws.Range("H31").Value = 38851.035
$ws.Range("I31").Value = 72897.5
$ws.Range("K31").Value = 72897.5
$ws.Range("M31").Value = -72602.5

$ws.Range("H34").Value = 38851.035
$ws.Range("I34").Value = 72897.5
$ws.Range("K34").Value = 72897.5
$ws.Range("M34").Value = -72695.5

$ws.Range("H43").Value = 37340.5
$ws.Range("J43").Value = 37340.5
$ws.Range("L43").Value = 37340.5
$ws.Range("N43").Value = -37708.5

$ws.Range("H94").Value = 2297.4167
$ws.Range("I94").Value = 2208.25
$ws.Range("J94").Value = 2342
$ws.Range("K94").Value = 2208.25
$ws.Range("L94").Value = 2342
$ws.Range("M94").Value = -1757.25
$ws.Range("N94").Value = -3244

$ws.Range("H95").Value = 17424.6
$ws.Range("J95").Value = 17424.6
$ws.Range("L95").Value = 17424.6
$ws.Range("N95").Value = -22916.6

$ws.Range("H101").Value = 37340.5
$ws.Range("J101").Value = 37340.5
$ws.Range("L101").Value = 37340.5
$ws.Range("N101").Value = -43830.5

$ws.Range("H105").Value = 1172.2858
$ws.Range("I105").Value = 1254.2222
$ws.Range("K105").Value = 1254.2222
$ws.Range("M105").Value = 492.7778000000001

$ws.Range("H132").Value = 2933.3438
$ws.Range("I132").Value = 2497.5173
$ws.Range("K132").Value = 7492.5519
$ws.Range("M132").Value = -4962.5519

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 518.125
$ws.Range("I86").Value = 515.8333
$ws.Range("K86").Value = 1547.4999
$ws.Range("M86").Value = -361.4999

$ws.Range("H89").Value = 518.125
$ws.Range("I89").Value = 515.8333
$ws.Range("K89").Value = 4642.4997
$ws.Range("M89").Value = 1285.5003

$ws.Range("H120").Value = 3499.5
$ws.Range("I120").Value = 3499.5
$ws.Range("K120").Value = 10498.5
$ws.Range("M120").Value = -5660.5

$ws.Range("H139").Value = 2562.375
$ws.Range("I139").Value = 2499.8572
$ws.Range("K139").Value = 7499.571599999999
$ws.Range("M139").Value = -2359.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 23386154
$ws.Range("I11").Value = 26728682
$ws.Range("J11").Value = 5002250
$ws.Range("K11").Value = 26728682
$ws.Range("L11").Value = 5002250
$ws.Range("M11").Value = -26728543
$ws.Range("N11").Value = -5002528

$ws.Range("H12").Value = 29250.25
$ws.Range("I12").Value = 37333.668
$ws.Range("K12").Value = 37333.668
$ws.Range("M12").Value = -37193.668

$ws.Range("H57").Value = 8312.5
$ws.Range("J57").Value = 14500
$ws.Range("L57").Value = 14500
$ws.Range("N57").Value = -16140

$ws.Range("H80").Value = 3683
$ws.Range("J80").Value = 3732.6667
$ws.Range("L80").Value = 3732.6667
$ws.Range("N80").Value = -5728.6667

$ws.Range("H83").Value = 3683
$ws.Range("J83").Value = 3732.6667
$ws.Range("L83").Value = 18663.3335
$ws.Range("N83").Value = -28647.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2756.611
$ws.Range("I22").Value = 1543.8
$ws.Range("J22").Value = 3223.077
$ws.Range("K22").Value = 1543.8
$ws.Range("L22").Value = 3223.077
$ws.Range("M22").Value = -1248.8
$ws.Range("N22").Value = -3813.077

$ws.Range("H27").Value = 2756.611
$ws.Range("I27").Value = 1543.8
$ws.Range("J27").Value = 3223.077
$ws.Range("K27").Value = 1543.8
$ws.Range("L27").Value = 3223.077
$ws.Range("M27").Value = -1436.8
$ws.Range("N27").Value = -3437.077

$ws.Range("H30").Value = 288.85715
$ws.Range("I30").Value = 288.85715
$ws.Range("K30").Value = 288.85715
$ws.Range("M30").Value = -180.85715

$ws.Range("H68").Value = 2938.889
$ws.Range("I68").Value = 2938.889
$ws.Range("K68").Value = 2938.889
$ws.Range("M68").Value = -2189.889

$ws.Range("H71").Value = 2938.889
$ws.Range("I71").Value = 2938.889
$ws.Range("K71").Value = 14694.445
$ws.Range("M71").Value = -10950.445

$ws.Range("H137").Value = 103000
$ws.Range("I137").Value = 90000
$ws.Range("J137").Value = 107333.336
$ws.Range("K137").Value = 90000
$ws.Range("L137").Value = 107333.336
$ws.Range("M137").Value = -84900
$ws.Range("N137").Value = -117533.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 29466.166
$ws.Range("I9").Value = 38200
$ws.Range("J9").Value = 11998.5
$ws.Range("K9").Value = 38200
$ws.Range("L9").Value = 11998.5
$ws.Range("M9").Value = -38060
$ws.Range("N9").Value = -12278.5

$ws.Range("H136").Value = 2393.111
$ws.Range("I136").Value = 1520.9
$ws.Range("K136").Value = 4562.700000000001
$ws.Range("M136").Value = -2012.700000000001
